$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "ID Competição" (column B) values from 64 to 264 for all data rows (B2:B209)
$ws.Range("B2:B209").Value = 264
